# Update the "Förändrad" (Changed) date column (C) for data rows 2-16
# from serial date 45207 (2023-10-08) to 45208 (2023-10-09).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 16; $row++) {
    $ws.Cells.Item($row, 3).Value = 45208
}
